# Rename the "preparation_temperature list" and "storage_temperature list"
# sheets, refresh their contents, update the related header comments and
# the data-validation rules that point at those lists.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the two list sheets.
# ---------------------------------------------------------------------
$prepSheet = $wb.Worksheets.Item("preparation_temperature list")
$prepSheet.Name = "preparation_condition list"

$storageSheet = $wb.Worksheets.Item("storage_temperature list")
$storageSheet.Name = "storage_method list"

# ---------------------------------------------------------------------
# 2. Replace the preparation_condition list values (7 rows).
# ---------------------------------------------------------------------
$prepSheet.Cells.Item(1, 1).Value = "frozen in liquid nitrogen"
$prepSheet.Cells.Item(2, 1).Value = "frozen in liquid nitrogen vapor"
$prepSheet.Cells.Item(3, 1).Value = "frozen in ice"
$prepSheet.Cells.Item(4, 1).Value = "frozen in dry ice"
$prepSheet.Cells.Item(5, 1).Value = "frozen at -20 C"
$prepSheet.Cells.Item(6, 1).Value = "ambient temperature"
$prepSheet.Cells.Item(7, 1).Value = "unknown"
$prepSheet.Rows.Item(8).Delete()

# ---------------------------------------------------------------------
# 3. Replace the storage_method list values (11 rows).
# ---------------------------------------------------------------------
$storageSheet.Cells.Item(1, 1).Value = "frozen in liquid nitrogen"
$storageSheet.Cells.Item(2, 1).Value = "frozen in liquid nitrogen vapor"
$storageSheet.Cells.Item(3, 1).Value = "frozen in ice"
$storageSheet.Cells.Item(4, 1).Value = "frozen in dry ice"
$storageSheet.Cells.Item(5, 1).Value = "frozen at -80 C"
$storageSheet.Cells.Item(6, 1).Value = "frozen at -20 C"
$storageSheet.Cells.Item(7, 1).Value = "refrigerator"
$storageSheet.Cells.Item(8, 1).Value = "ambient temperature"
$storageSheet.Cells.Item(9, 1).Value = "incubated at 37 C"
$storageSheet.Cells.Item(10, 1).Value = "none"
$storageSheet.Cells.Item(11, 1).Value = "unknown"
$storageSheet.Rows.Item(12).Delete()

# ---------------------------------------------------------------------
# 4. Update the header labels and comments on the "Export as TSV" sheet.
# ---------------------------------------------------------------------
$tsv = $wb.Worksheets.Item("Export as TSV")
$tsv.Cells.Item(1, 7).Value = "preparation_condition"
$tsv.Cells.Item(1, 11).Value = "storage_method"

$tsv.Range("G1").Comment.Text("The condition under which the preparation occurred, such as whether the sample was placed in dry ice during the preparation.")
$tsv.Range("K1").Comment.Text("The method by which the sample was stored, after preparation and before the assay was performed.")

# ---------------------------------------------------------------------
# 5. Point the data validations at the renamed lists / updated ranges.
# ---------------------------------------------------------------------
$gRange = $tsv.Range("G2:G1048576")
$gRange.Validation.Delete()
$gRange.Validation.Add(3, 1, 1, "='preparation_condition list'!`$A`$1:`$A`$7")
$gRange.Validation.ErrorTitle = "Value must come from list"
$gRange.Validation.ErrorMessage = "Value must come from preparation_condition list."
$gRange.Validation.IgnoreBlank = $true
$gRange.Validation.InCellDropdown = $true
$gRange.Validation.ShowInput = $true
$gRange.Validation.ShowError = $true

$kRange = $tsv.Range("K2:K1048576")
$kRange.Validation.Delete()
$kRange.Validation.Add(3, 1, 1, "='storage_method list'!`$A`$1:`$A`$11")
$kRange.Validation.ErrorTitle = "Value must come from list"
$kRange.Validation.ErrorMessage = "Value must come from storage_method list."
$kRange.Validation.IgnoreBlank = $true
$kRange.Validation.InCellDropdown = $true
$kRange.Validation.ShowInput = $true
$kRange.Validation.ShowError = $true

Write-Host "Done"
